$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 295, shifting existing rows 295-329 down to 296-330.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new data entry.
$ws.Cells.Item(295, 1).Value2 = 10
$ws.Cells.Item(295, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(295, 3).Value2 = "La Araucanía"
$ws.Cells.Item(295, 4).Value2 = 44995
$ws.Cells.Item(295, 5).Value2 = 9
$ws.Cells.Item(295, 6).Value2 = "Fruta"
$ws.Cells.Item(295, 7).Value2 = 100103
$ws.Cells.Item(295, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(295, 9).Value2 = 100103002
$ws.Cells.Item(295, 10).Value2 = "Ciruela"
$ws.Cells.Item(295, 11).Value2 = "Blue Giant"
$ws.Cells.Item(295, 12).Value2 = "Primera"
$ws.Cells.Item(295, 13).Value2 = 200
$ws.Cells.Item(295, 14).Value2 = 14000
$ws.Cells.Item(295, 15).Value2 = 15000
$ws.Cells.Item(295, 16).Value2 = 14500
$ws.Cells.Item(295, 17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(295, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(295, 19).Value2 = 806
$ws.Cells.Item(295, 20).Value2 = 18
